$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy formatting from row 14 to row 15 first, so the new row's styles
# (s="2", s="3") already match before any values are written.
$ws.Range("A14:E14").Copy() | Out-Null
$ws.Range("A15:E15").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Add new row 15: ID 14, Date 2016-12-16 (serial 42720), Title (new shared string),
# Description blank, Status COMPLETED
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = 42720
$ws.Range("C15").Value = "Finish threading Select from CmdLine to DataStore"
$ws.Range("E15").Value = "COMPLETED"

# Update selection to match target diff
$ws.Range("C22").Select() | Out-Null
